# faturamento_diario.xlsx update
#
# 1) Correct the 19/07/2025 total_venda value (row 20, column B).
# 2) Insert two new daily rows for July (days 21 and 22) right after the
#    existing July block (row 21), pushing every row below it down by 2 -
#    this naturally re-aligns the June/May/April blocks to their new row
#    numbers without touching their contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix existing value -----------------------------------------------
$ws.Range("B20").Value = 15007.65

# --- insert two fresh rows for July 21 and July 22 ---------------------
$ws.Rows("22:23").Insert()

$periodoJul = $ws.Range("E2").Value()  # reuse the existing "07/2025" shared string

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 26242.62
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 2025
$ws.Range("E22").Value = $periodoJul

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 60
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 2025
$ws.Range("E23").Value = $periodoJul
